# "make import excel for kepesertaan"
# Build the Kepesertaan (participant) import header table above the
# existing single data row, matching the target layout:
#   rows 1-6   : blank spacer rows (pre-formatted)
#   row 2,I    : title line (shared string)
#   rows 7-9   : merged table header ("NO", "NO Peserta", "NIP", ...)
#   row 10     : original data row (shifted down from row 1)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------
# 1. Push the existing single data row down from row 1 to row 10,
#    inserting 9 blank rows above it. Existing styles/values on
#    that row travel with it automatically.
# ---------------------------------------------------------------
$ws.Rows("1:9").Insert()
$ws.Range("A10").ClearFormats()

# ---------------------------------------------------------------
# 2. Formatting "source" cells already carrying the fonts we need
#    (re-use them via Copy + PasteSpecial(Formats) so the engine
#    reuses existing font/xf entries instead of inventing new ones
#    with spurious theme colors).
# ---------------------------------------------------------------
# B10 = Arial 10, no border, no alignment override (style used by columns B:E)

# Rows 1-6: spacer rows above the header table
$ws.Range("A1:E1").Value = $null
$ws.Range("B10").Copy()
$ws.Range("A1:E6").PasteSpecial(-4122)

$ws.Range("B10").Copy()
$ws.Range("F1:I6").PasteSpecial(-4122)
$ws.Range("F1:I6").HorizontalAlignment = -4108

$ws.Range("B10").Copy()
$ws.Range("J1:L6").PasteSpecial(-4122)
$ws.Range("J1:L6").NumberFormat = "#,##0"

# Title text in I2 (document/report header line)
$ws.Range("I2").Value = "PAHUPL-PUSAT800-201812, PAHUPL, 01-12-2018, ISMAIL, 17-12-2018 11:51:56"

# ---------------------------------------------------------------
# 3. Header table, rows 7-9 (merged column headers with borders)
# ---------------------------------------------------------------

# --- "NO" column: B7:B9 merged, full box border ---
$ws.Range("B10").Copy()
$ws.Range("B7:B9").PasteSpecial(-4122)
$ws.Range("B7:B9").HorizontalAlignment = -4108
$ws.Range("B7:B9").VerticalAlignment = -4108
$ws.Range("B7:B9").Borders.Item(7).LineStyle = 1
$ws.Range("B7:B9").Borders.Item(8).LineStyle = 1
$ws.Range("B7:B9").Borders.Item(9).LineStyle = 1
$ws.Range("B7:B9").Borders.Item(10).LineStyle = 1
$ws.Range("B7").Value = "NO"
$ws.Range("B7:B9").Merge()

# --- "NIP" column header top cell style (border, no bottom) reused for C7 & H7 ---
$ws.Range("B10").Copy()
$ws.Range("C7").PasteSpecial(-4122)
$ws.Range("C7").HorizontalAlignment = -4108
$ws.Range("C7").VerticalAlignment = -4108
$ws.Range("C7").Borders.Item(7).LineStyle = 1
$ws.Range("C7").Borders.Item(10).LineStyle = 1
$ws.Range("C7").Borders.Item(8).LineStyle = 1
$ws.Range("C7").Value = "NO Peserta"

$ws.Range("C7").Copy()
$ws.Range("H7").PasteSpecial(-4122)
$ws.Range("H7").Value = "Golongan"

# --- "NIP" column D7:D9 (same full-box style as B) ---
$ws.Range("B7").Copy()
$ws.Range("D7:D9").PasteSpecial(-4122)
$ws.Range("D7").Value = "NIP"
$ws.Range("D7:D9").Merge()

$ws.Range("B7").Copy()
$ws.Range("E7:E9").PasteSpecial(-4122)
$ws.Range("E7").Value = "Nama Peserta"
$ws.Range("E7:E9").Merge()

# --- "Unit kerja" spanning F7:G7 (full box, horizontal-only center) ---
$ws.Range("B10").Copy()
$ws.Range("F7:G7").PasteSpecial(-4122)
$ws.Range("F7:G7").HorizontalAlignment = -4108
$ws.Range("F7:G7").Borders.Item(7).LineStyle = 1
$ws.Range("F7:G7").Borders.Item(10).LineStyle = 1
$ws.Range("F7:G7").Borders.Item(8).LineStyle = 1
$ws.Range("F7:G7").Borders.Item(9).LineStyle = 1
$ws.Range("F7").Value = "Unit kerja"
$ws.Range("F7:G7").Merge()

# --- "Stt" I7 (full box, horizontal-only center; single cell, not merged) ---
$ws.Range("F7").Copy()
$ws.Range("I7").PasteSpecial(-4122)
$ws.Range("I7").Value = "Stt"

# --- "Gaji pokok" / "Gaji Pns" J7:J9 / K7:K9 (full box, numFmt, centered) ---
$ws.Range("B10").Copy()
$ws.Range("J7:J9").PasteSpecial(-4122)
$ws.Range("J7:J9").HorizontalAlignment = -4108
$ws.Range("J7:J9").VerticalAlignment = -4108
$ws.Range("J7:J9").NumberFormat = "#,##0"
$ws.Range("J7:J9").Borders.Item(7).LineStyle = 1
$ws.Range("J7:J9").Borders.Item(10).LineStyle = 1
$ws.Range("J7:J9").Borders.Item(8).LineStyle = 1
$ws.Range("J7:J9").Borders.Item(9).LineStyle = 1
$ws.Range("J7").Value = "Gaji pokok"
$ws.Range("J7:J9").Merge()

$ws.Range("J7").Copy()
$ws.Range("K7:K9").PasteSpecial(-4122)
$ws.Range("K7").Value = "Gaji Pns"
$ws.Range("K7:K9").Merge()

# --- "PHDP" L7:L9 (border without bottom on the top cell, numFmt, centered) ---
$ws.Range("B10").Copy()
$ws.Range("L7").PasteSpecial(-4122)
$ws.Range("L7").HorizontalAlignment = -4108
$ws.Range("L7").VerticalAlignment = -4108
$ws.Range("L7").NumberFormat = "#,##0"
$ws.Range("L7").Borders.Item(7).LineStyle = 1
$ws.Range("L7").Borders.Item(10).LineStyle = 1
$ws.Range("L7").Borders.Item(8).LineStyle = 1
$ws.Range("L7").Value = "PHDP"

# ---------------------------------------------------------------
# Row 8: "Saat ini" / "Mutasi" sub-headers
# ---------------------------------------------------------------
# F8:F9 merged, full box border (re-uses B-style)
$ws.Range("B7").Copy()
$ws.Range("F8:F9").PasteSpecial(-4122)
$ws.Range("F8").Value = "Saat ini"
$ws.Range("F8:F9").Merge()

# G8 "Mutasi " - border without bottom (same pattern as C7/H7)
$ws.Range("C7").Copy()
$ws.Range("G8").PasteSpecial(-4122)
$ws.Range("G8").Value = "Mutasi "

# C8 / H8 - middle cell of vertical merge: sides only, no top, no bottom
$ws.Range("B10").Copy()
$ws.Range("C8").PasteSpecial(-4122)
$ws.Range("C8").HorizontalAlignment = -4108
$ws.Range("C8").VerticalAlignment = -4108
$ws.Range("C8").Borders.Item(7).LineStyle = 1
$ws.Range("C8").Borders.Item(10).LineStyle = 1

$ws.Range("C8").Copy()
$ws.Range("H8").PasteSpecial(-4122)

# I8:I9 merged, full box (re-use B-style)
$ws.Range("B7").Copy()
$ws.Range("I8:I9").PasteSpecial(-4122)
$ws.Range("I8").Value = "Kawin"
$ws.Range("I8:I9").Merge()

# L8 - middle cell (sides only, no top/bottom) with numFmt
$ws.Range("C8").Copy()
$ws.Range("L8").PasteSpecial(-4122)
$ws.Range("L8").NumberFormat = "#,##0"

# ---------------------------------------------------------------
# Row 9: "Dari" sub-header + bottom border cells
# ---------------------------------------------------------------
# G9 "Dari" - border with bottom, no top (mirrors C9/H9 pattern)
$ws.Range("B10").Copy()
$ws.Range("C9").PasteSpecial(-4122)
$ws.Range("C9").HorizontalAlignment = -4108
$ws.Range("C9").VerticalAlignment = -4108
$ws.Range("C9").Borders.Item(7).LineStyle = 1
$ws.Range("C9").Borders.Item(10).LineStyle = 1
$ws.Range("C9").Borders.Item(9).LineStyle = 1

$ws.Range("C9").Copy()
$ws.Range("H9").PasteSpecial(-4122)

$ws.Range("C9").Copy()
$ws.Range("G9").PasteSpecial(-4122)
$ws.Range("G9").HorizontalAlignment = -4108
$ws.Range("G9").Value = "Dari"

# L9 bottom cell, with numFmt
$ws.Range("C9").Copy()
$ws.Range("L9").PasteSpecial(-4122)
$ws.Range("L9").NumberFormat = "#,##0"

# ---------------------------------------------------------------
# Final sheet view / selection state
# ---------------------------------------------------------------
$ws.Range("K7:K9").Select()
